$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.495.92"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.811.41"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.10%  "
$c = $ws.Range("D5")
$c.Value = "'225.55"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("E7").Value = "  -0.07%  "
$c = $ws.Range("D8")
$c.Value = "'38.28"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +6.22%  "
$ws.Range("E9").Value = "  -4.01%  "
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "2.073.04"
$ws.Range("E12").Value = "  +0.53%  "
$c = $ws.Range("D13")
$c.Value = "'11.24"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "1.809.71"
$ws.Range("E14").Value = "  +0.74%  "
$c = $ws.Range("D15")
$c.Value = "'0.633"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").Value = "34.444.04"
$ws.Range("E16").Value = "  +0.36%  "
$c = $ws.Range("D17")
$c.Value = "'4.43"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.51%  "
$c = $ws.Range("D18")
$c.Value = "'68.36"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$c = $ws.Range("D19")
$c.Value = "'243.16"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("E22").Value = "  -0.05%  "
$c = $ws.Range("D23")
$c.Value = "'4.13"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("E24").Value = "  +3.53%  "
$c = $ws.Range("D25")
$c.Value = "'170.44"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.01%  "
$c = $ws.Range("D26")
$c.Value = "'7.81"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.07%  "
$c = $ws.Range("D27")
$c.Value = "'17.61"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.48%  "
$c = $ws.Range("D28")
$c.Value = "'0.121"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D30")
$c.Value = "'3.80"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D31")
$c.Value = "'1.23"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.40%  "
$c = $ws.Range("D32")
$c.Value = "'0.0518"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.62%  "
$c = $ws.Range("D33")
$c.Value = "'3.87"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -4.48%  "
$c = $ws.Range("D34")
$c.Value = "'1.83"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").Value = "1.361.33"
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("E36").Value = "  -4.52%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D38")
$c.Value = "'0.0187"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D39")
$c.Value = "'2.34"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -5.20%  "
$c = $ws.Range("D40")
$c.Value = "'2.45"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D41")
$c.Value = "'0.953"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D42")
$c.Value = "'81.91"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D43")
$c.Value = "'1.21"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.42%  "
$c = $ws.Range("D44")
$c.Value = "'2.81"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.31%  "
$c = $ws.Range("D45")
$c.Value = "'13.82"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("D47").Value = "1.974.22"
$ws.Range("E47").Value = "  +0.57%  "
$c = $ws.Range("D48")
$c.Value = "'5.78"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("E49").Value = "  -0.06%  "
$c = $ws.Range("D50")
$c.Value = "'102.46"
$c.Style = "Normal"
$ws.Range("D51").Value = "0.0₆0121"
$ws.Range("E51").Value = "  -5.10%  "
